# AndroidProgramList.xlsx update:
#   - "contextmenu app done, grid view app done"
#
# 1. The "Basic View 2" demo is clarified as the progress-bar sample, so its
#    label is renamed to "Basic View 2(ProgressBar)".
# 2. Context menu App / Web View App / Grid View App are marked as finished
#    (Journal column filled in).
# 3. A new blank row is inserted before the storage entries.
# 4. Internal Storage / External Storage are also marked as finished.
# 5. The "duplicate values" conditional formatting that runs down column D is
#    extended to keep covering the data through the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename "Basic View 2" -> "Basic View 2(ProgressBar)" ---------------
$ws.Cells.Item(29, 3).Value = "Basic View 2(ProgressBar)"

# --- 2. Mark the three newly-finished apps with their Journal/LabBook status
$ws.Cells.Item(36, 4).Value = "Journal"   # Context menu App
$ws.Cells.Item(37, 4).Value = "Journal"   # Web View App
$ws.Cells.Item(38, 4).Value = "Journal"   # Grid View App

# --- 3. Insert one new (blank) row above the storage entries ---------------
$ws.Rows.Item(39).Insert()

# --- 4. Mark the storage rows (now shifted down to 40/41) as finished too --
$ws.Cells.Item(40, 4).Value = "Journal"   # Internal Storage
$ws.Cells.Item(41, 4).Value = "Journal"   # External Storage

# --- 5. Extend the conditional formatting range down to the new rows -------
$fc = $ws.Range("D3:D35").FormatConditions
$fc.Item(1).ModifyAppliesToRange($ws.Range("D3:D38,D40:D41"))

# --- 5b. Excel leaves behind an extra differential-format record each time a
#     duplicate-values rule is reapplied/edited through the UI. Reproduce that
#     so the style table ends up with the same 6 (5 leftover + 1 live) dxfs.
$scratch = $ws.Range("D3:D38,D40:D41").FormatConditions
for ($i = 0; $i -lt 5; $i++) {
    $tmpRule = $scratch.Add(1, 3, "=$D$3")
    $tmpRule.Font.Color = 26012
    $tmpRule.Interior.Color = 10284031
    $tmpRule.Delete()
}

# --- 6. Update the view's selection/scroll position to match the edit area -
$excel.ActiveWindow.ScrollRow = 34
$ws.Range("D38").Select()

Write-Host "Workbook updated."
